# Updated cryptos list on Thu Nov 23 07:56:01 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "37.374.56"
$ws.Range("E2").Value = "  +2.48%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.065.49"
$ws.Range("E3").Value = "  +3.66%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.00%  "

# Row 5 - BNB
$ws.Range("D5").Value = "234.77"

# Row 6 - XRP
$ws.Range("D6").Value = "0.615"
$ws.Range("E6").Value = "  +2.63%  "

# Row 7 - Solana
$ws.Range("D7").Value = "'57.90"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.51%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.01%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "0.382"
$ws.Range("E9").Value = "  +3.14%  "

# Row 10 - OKB
$ws.Range("D10").Value = "58.91"
$ws.Range("E10").Value = "  +1.87%  "

# Row 11 - Dogecoin
$ws.Range("D11").Value = "'0.0760"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.85%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  +2.93%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "2.369.50"
$ws.Range("E13").Value = "  +3.47%  "

# Row 14 - Chainlink
$ws.Range("E14").Value = "  +2.42%  "

# Row 15 - Avalanche
$ws.Range("D15").Value = "21.09"
$ws.Range("E15").Value = "  +4.19%  "

# Row 16 - Polygon
$ws.Range("D16").Value = "0.777"
$ws.Range("E16").Value = "  +2.57%  "

# Row 17 - Polkadot
$ws.Range("D17").Value = "5.18"
$ws.Range("E17").Value = "  +2.41%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "2.066.04"
$ws.Range("E18").Value = "  +2.83%  "

# Row 19 - WrappedBTC
$ws.Range("D19").Value = "37.579.29"
$ws.Range("E19").Value = "  +3.10%  "

# Row 20 - Uniswap
$ws.Range("D20").Value = "6.16"
$ws.Range("E20").Value = "  +16.67%  "

# Row 21 - Litecoin
$ws.Range("D21").Value = "70.35"
$ws.Range("E21").Value = "  +3.90%  "

# Row 22 - ShibaInu
$ws.Range("E22").Value = "  +1.42%  "

# Row 23 - BitcoinCash
$ws.Range("D23").Value = "226.86"
$ws.Range("E23").Value = "  +2.31%  "

# Row 24 - Dai
$ws.Range("E24").Value = "  -0.09%  "

# Row 25 - PancakeSwap
$ws.Range("D25").Value = "2.44"
$ws.Range("E25").Value = "  +2.17%  "

# Row 26 - Toncoin
$ws.Range("E26").Value = "  +0.99%  "

# Row 27 - Monero
$ws.Range("D27").Value = "165.41"
$ws.Range("E27").Value = "  +2.20%  "

# Row 28 - ImmutableX
$ws.Range("E28").Value = "  +12.59%  "

# Row 29 - Cosmos
$ws.Range("E29").Value = "  +2.57%  "

# Row 30 - was EthereumClassic, now Kaspa
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").Value = "0.128"
$ws.Range("E30").Value = "  +1.50%  "

# Row 31 - was Kaspa, now EthereumClassic
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "19.15"
$ws.Range("E31").Value = "  +2.06%  "

# Row 32 - Stellar
$ws.Range("D32").Value = "0.119"
$ws.Range("E32").Value = "  +2.19%  "

# Row 33 - Filecoin
$ws.Range("D33").Value = "'4.50"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.85%  "

# Row 34 - Hedera
$ws.Range("E34").Value = "  +2.96%  "

# Row 35 - LidoDAOToken
$ws.Range("D35").Value = "2.56"
$ws.Range("E35").Value = "  +9.50%  "

# Row 36 - InternetComputer(DFINITY)
$ws.Range("D36").Value = "4.55"
$ws.Range("E36").Value = "  +7.22%  "

# Row 37 - was BinanceUSD, now RenderToken
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").Value = "3.37"
$ws.Range("E37").Value = "  +0.57%  "

# Row 38 - was RenderToken, now BinanceUSD
$ws.Range("B38").Value = "BinanceUSD"
$ws.Range("C38").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D38").Value = "'1.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.00%  "

# Row 39 - WEMIXToken
$ws.Range("D39").Value = "1.79"
$ws.Range("E39").Value = "  +1.08%  "

# Row 40 - THORChain
$ws.Range("E40").Value = "  +4.12%  "

# Row 41 - HuobiToken
$ws.Range("E41").Value = "  -1.53%  "

# Row 42 - Cronos
$ws.Range("D42").Value = "0.0967"
$ws.Range("E42").Value = "  +3.36%  "

# Row 43 - FTXToken
$ws.Range("D43").Value = "4.39"
$ws.Range("E43").Value = "  +20.77%  "

# Row 44 - Aave
$ws.Range("D44").Value = "96.03"
$ws.Range("E44").Value = "  +8.01%  "

# Row 45 - Maker
$ws.Range("D45").Value = "1.453.46"
$ws.Range("E45").Value = "  -0.28%  "

# Row 46 - VeChain
$ws.Range("E46").Value = "  +4.26%  "

# Row 47 - TrustWalletToken
$ws.Range("E47").Value = "  +6.07%  "

# Row 48 - InjectiveProtocol
$ws.Range("E48").Value = "  +4.44%  "

# Row 49 - ARBITRUM
$ws.Range("E49").Value = "  +4.22%  "

# Row 50 - FraxShare
$ws.Range("E50").Value = "  +6.81%  "

# Row 51 - MXToken
$ws.Range("E51").Value = "  +1.93%  "
